$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'65.862.68"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  +0.14%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'2.661.34"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  -0.68%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("E4").Value = "'  +0.07%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'598.68"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  -0.34%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'159.26"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  +1.45%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'0.644"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'  +3.74%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("D9").Value = "'0.127"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'  -3.08%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("E10").Value = "'  +0.21%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'5.87"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'  -0.25%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("E12").Value = "'  +1.51%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'29.08"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  -1.05%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'0.0000196"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  -1.92%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'3.145.85"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  -0.44%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'65.806.68"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  +0.28%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'2.648.94"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  -0.59%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'12.63"
$ws.Range("D18").Style = "Normal"
$ws.Range("D19").Value = "'4.82"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  +0.08%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'352.89"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  +0.04%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'7.49"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  -1.53%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("E22").Value = "'  -0.11%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'70.13"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  +0.40%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'1.81"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  +10.04%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'0.0000113"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  +1.78%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'9.66"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'  -0.07%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = "'1.62"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'  +1.74%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("D28").Value = "'577.16"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'  +8.92%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("D29").Value = "'8.16"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  +1.37%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("E30").Value = "'  -1.86%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("E31").Value = "'  -0.21%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("E32").Value = "'  +0.79%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("D33").Value = "'1.81"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "'  +2.29%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("D34").Value = "'6.74"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "'  +3.84%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("D35").Value = "'5.54"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  +0.54%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("D36").Value = "'0.423"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "'  -0.23%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("D37").Value = "'20.61"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "'  -0.45%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("E38").Value = "'  +0.00%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("D39").Value = "'1.96"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'  +0.96%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'155.04"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'  -1.89%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("E41").Value = "'  +8.56%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'161.92"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  -1.52%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'4.11"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  -0.93%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'0.0618"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  +0.93%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'23.53"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  +2.73%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'0.644"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  +0.20%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("E47").Value = "'  +0.45%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("E48").Value = "'  +1.12%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = "'19.85"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  -1.81%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = "'0.0₆0247"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  -6.91%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = "'0.818"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  +0.19%  "
$ws.Range("E51").Style = "Normal"
